$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records are inserted at the top of the data block
# (rows 519-520), pushing the existing rows 519-536 down to 521-538.
$ws.Rows.Item(519).Insert()
$ws.Rows.Item(519).Insert()

# Row 519 - new record
$ws.Range("A519").Value() = 4
$ws.Range("B519").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C519").Value() = "Los Lagos"
$ws.Range("D519").Value() = 44747
$ws.Range("E519").Value() = 10
$ws.Range("F519").Value() = "Fruta"
$ws.Range("G519").Value() = 100108
$ws.Range("H519").Value() = "Tropicales y subtropicales"
$ws.Range("I519").Value() = 100108006
$ws.Range("J519").Value() = "Plátano"
$ws.Range("K519").Value() = "Sin especificar"
$ws.Range("L519").Value() = "Pintón"
$ws.Range("M519").Value() = 500
$ws.Range("N519").Value() = 25000
$ws.Range("O519").Value() = 25000
$ws.Range("P519").Value() = 25000
$ws.Range("Q519").Value() = "$/caja 20 kilos"
$ws.Range("R519").Value() = "Ecuador"
$ws.Range("S519").Value() = 1250
$ws.Range("T519").Value() = 20

# Row 520 - new record
$ws.Range("A520").Value() = 4
$ws.Range("B520").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C520").Value() = "Los Lagos"
$ws.Range("D520").Value() = 44747
$ws.Range("E520").Value() = 10
$ws.Range("F520").Value() = "Fruta"
$ws.Range("G520").Value() = 100108
$ws.Range("H520").Value() = "Tropicales y subtropicales"
$ws.Range("I520").Value() = 100108006
$ws.Range("J520").Value() = "Plátano"
$ws.Range("K520").Value() = "Sin especificar"
$ws.Range("L520").Value() = "Primera Pintón"
$ws.Range("M520").Value() = 1000
$ws.Range("N520").Value() = 26000
$ws.Range("O520").Value() = 27000
$ws.Range("P520").Value() = 26500
$ws.Range("Q520").Value() = "$/caja 20 kilos"
$ws.Range("R520").Value() = "Ecuador"
$ws.Range("S520").Value() = 1325
$ws.Range("T520").Value() = 20

Write-Output "applied"
